# Re-commit the three ability-changelog "effect" strings so that they sort
# after the rest of the shared-string pool (i.e. "move to the end of the
# table"), matching the source data regeneration that produced the target
# workbook. The three affected rows are:
#   D12 -> "Inflicts only 1/16 of the attacker's maximum [HP]{mechanic:hp} in damage."
#   D24 -> "Does not affect friendly Pokémon's moves that target all other Pokémon.  This ability's presence is not announced upon entering battle."
#   D28 -> "Doubles []{move:cut}'s grass-cutting radius on the overworld if any party Pokémon has this ability."
# No other cell text changes - this is purely a re-commit of existing data
# (matches the commit message "Commit data, update entity, service, repository").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$inflicts = "Inflicts only 1/16 of the attacker's maximum [HP]{mechanic:hp} in damage."
$friendly = "Does not affect friendly Pokémon's moves that target all other Pokémon.  This ability's presence is not announced upon entering battle."
$doubles  = "Doubles []{move:cut}'s grass-cutting radius on the overworld if any party Pokémon has this ability."

# Re-write in the same relative order the rows appear in the refreshed
# shared-string table so the data is re-committed identically to source.
$ws.Range("D12").Value = $inflicts
$ws.Range("D24").Value = $friendly
$ws.Range("D28").Value = $doubles
